$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 125 -- this shifts the existing rows 125-150
# down to 126-151, matching the dimension growth from A1:T150 to A1:T151.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly price record.
$ws.Cells.Item(125, 1).Value = 10
$ws.Cells.Item(125, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(125, 3).Value = "La Araucanía"
$ws.Cells.Item(125, 4).Value = 45258
$ws.Cells.Item(125, 5).Value = 9
$ws.Cells.Item(125, 6).Value = "Fruta"
$ws.Cells.Item(125, 7).Value = 100101
$ws.Cells.Item(125, 8).Value = "Berries"
$ws.Cells.Item(125, 9).Value = 100101001
$ws.Cells.Item(125, 10).Value = "Arándano (blue)"
$ws.Cells.Item(125, 11).Value = "Sin especificar"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 200
$ws.Cells.Item(125, 14).Value = 4000
$ws.Cells.Item(125, 15).Value = 4000
$ws.Cells.Item(125, 16).Value = 4000
$ws.Cells.Item(125, 17).Value = "$/kilo"
$ws.Cells.Item(125, 18).Value = "Región del Maule"
$ws.Cells.Item(125, 19).Value = 4000
$ws.Cells.Item(125, 20).Value = 1
